$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new EVF ticker row (row 33) with its bank type classification
$ws.Range("A33").Value = "EVF"
$ws.Range("B33").Value = "Private_3"

# Restore the view/selection state to match the saved workbook
$ws.Range("B33").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
